$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 673, pushing all existing rows (673..732) down to (675..734)
$ws.Range("A673:A674").EntireRow.Insert()

# Populate new row 673
$ws.Cells.Item(673, 1).Value = 7
$ws.Cells.Item(673, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(673, 3).Value = "Ñuble"
$ws.Cells.Item(673, 4).Value = 44769
$ws.Cells.Item(673, 5).Value = 16
$ws.Cells.Item(673, 6).Value = "Fruta"
$ws.Cells.Item(673, 7).Value = 100102
$ws.Cells.Item(673, 8).Value = "Cítricos"
$ws.Cells.Item(673, 9).Value = 100102003
$ws.Cells.Item(673, 10).Value = "Limón"
$ws.Cells.Item(673, 11).Value = "Sin especificar"
$ws.Cells.Item(673, 12).Value = "1a amarillo"
$ws.Cells.Item(673, 13).Value = 160
$ws.Cells.Item(673, 14).Value = 4500
$ws.Cells.Item(673, 15).Value = 5000
$ws.Cells.Item(673, 16).Value = 4750
$ws.Cells.Item(673, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(673, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(673, 19).Value = 297
$ws.Cells.Item(673, 20).Value = 16

# Populate new row 674
$ws.Cells.Item(674, 1).Value = 7
$ws.Cells.Item(674, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(674, 3).Value = "Ñuble"
$ws.Cells.Item(674, 4).Value = 44769
$ws.Cells.Item(674, 5).Value = 16
$ws.Cells.Item(674, 6).Value = "Fruta"
$ws.Cells.Item(674, 7).Value = 100102
$ws.Cells.Item(674, 8).Value = "Cítricos"
$ws.Cells.Item(674, 9).Value = 100102003
$ws.Cells.Item(674, 10).Value = "Limón"
$ws.Cells.Item(674, 11).Value = "Sin especificar"
$ws.Cells.Item(674, 12).Value = "2a amarillo"
$ws.Cells.Item(674, 13).Value = 80
$ws.Cells.Item(674, 14).Value = 4000
$ws.Cells.Item(674, 15).Value = 4000
$ws.Cells.Item(674, 16).Value = 4000
$ws.Cells.Item(674, 17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(674, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(674, 19).Value = 250
$ws.Cells.Item(674, 20).Value = 16
